# Applies the "Atualizado por script em 20-12-2023 02:45" update:
#  - rows 78/79 have their home/away match data (cols F:V) swapped
#  - rows 88/89 have their home/away match data (cols F:V) swapped
#  - three new match rows (171, 172, 173) are appended at the end

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($sheet, [int]$row1, [int]$row2)

    $cols = @("F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
    foreach ($col in $cols) {
        $addr1 = "$col$row1"
        $addr2 = "$col$row2"
        $v1 = $sheet.Range($addr1).Value2
        $v2 = $sheet.Range($addr2).Value2
        $sheet.Range($addr1).Value = $v2
        $sheet.Range($addr2).Value = $v1
    }
}

# --- Swap rows 78 and 79 ---
Swap-RowData $ws 78 79

# --- Swap rows 88 and 89 ---
Swap-RowData $ws 88 89

# --- Append new rows 171, 172, 173 ---
# Copy formatting from the last existing row (170) so style indices
# (bold/border on col A, datetime format on col E) match.
$ws.Range("A170:V170").Copy()
$ws.Range("A171:V173").PasteSpecial(-4122)

function Set-Row {
    param($sheet, [int]$r, [array]$values)

    $cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V")
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $sheet.Range("$($cols[$i])$r").Value = $values[$i]
    }
}

Set-Row $ws 171 @(
    170, "spain", "laliga", "2023-2024", 45279.79166666666,
    "Rayo Vallecano", 0, "Valencia", 1,
    2.24, "09/12/2023 22:02", 2.28, "19/12/2023 18:56",
    3.22, "09/12/2023 22:02", 3.17, "19/12/2023 18:58",
    3.25, "09/12/2023 22:02", 3.66, "19/12/2023 18:52",
    "https://www.betexplorer.com/football/spain/laliga/rayo-vallecano-valencia/QVVqd7km/"
)

Set-Row $ws 172 @(
    171, "spain", "laliga", "2023-2024", 45279.89583333334,
    "Atl. Madrid", 3, "Getafe", 3,
    1.31, "10/12/2023 00:01", 1.48, "19/12/2023 21:28",
    4.91, "10/12/2023 00:01", 4.54, "19/12/2023 21:29",
    9, "10/12/2023 00:01", 7.35, "19/12/2023 21:29",
    "https://www.betexplorer.com/football/spain/laliga/atl-madrid-getafe/EwmSwQZI/"
)

Set-Row $ws 173 @(
    172, "spain", "laliga", "2023-2024", 45279.89583333334,
    "Granada CF", 0, "Sevilla", 3,
    2.76, "10/12/2023 00:01", 3.12, "19/12/2023 21:27",
    3.33, "10/12/2023 00:01", 3.17, "19/12/2023 21:05",
    2.48, "10/12/2023 00:01", 2.55, "19/12/2023 21:25",
    "https://www.betexplorer.com/football/spain/laliga/granada-cf-sevilla/KhVmeR4g/"
)
